$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.002.67'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.83%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.830.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.27%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4655'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.42%  '

# Row 8
$ws.Range("E8").Value = '  +1.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07437'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8744'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.11%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.03'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.08%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07865'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.66%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.834.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.25%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.371'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.82%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.597'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.82%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.52%  '

# Row 17
$ws.Range("E17").Value = '  +0.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008963'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.98%  '

# Row 19
$ws.Range("E19").Value = '  -0.28%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.19%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.045.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.169'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.14%  '

# Row 23
$ws.Range("E23").Value = '  -0.24%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.063.26'
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.835'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.95%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.80%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.104'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.41%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.132'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.72%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.93%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08881'
$ws.Range("D31").Style = "Normal"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.971'
$ws.Range("D32").Style = "Normal"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7306'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.91%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.454'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.60%  '

# Row 35
$ws.Range("E35").Value = '  -3.09%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.499'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.58%  '

# Row 37
$ws.Range("E37").Value = '  -1.06%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01952'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.16%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05249'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.55%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.365'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.07%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.929'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.70%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5189'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.68%  '

# Row 43
$ws.Range("E43").Value = '  -1.35%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8574'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -15.17%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.233'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.79%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4860'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.82%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.75%  '

# Row 48
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.31%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.17%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.626'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.12%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06252'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.80%  '
